$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 180
$ws.Cells.Item(179, 1).Copy($ws.Cells.Item(180, 1))
$ws.Cells.Item(180, 1).Value = 178
$ws.Cells.Item(180, 2).Value = "Fiorentina"
$ws.Cells.Item(180, 3).Value = "Napoli"
$ws.Cells.Item(180, 4).Value = 0
$ws.Cells.Item(180, 5).Value = 3
$ws.Cells.Item(180, 6).Value = 1.19
$ws.Cells.Item(180, 7).Value = 1.64
$ws.Cells.Item(180, 8).Value = 1.47
$ws.Cells.Item(180, 9).Value = 1.84
$ws.Cells.Item(180, 10).Value = 1
$ws.Cells.Item(180, 11).Value = 1
$ws.Cells.Item(180, 12).Value = 0.28
$ws.Cells.Item(180, 13).Value = 0.2
$ws.Cells.Item(180, 14).Value = 0.48
$ws.Cells.Item(180, 15).Value = 3

# Row 181
$ws.Cells.Item(180, 1).Copy($ws.Cells.Item(181, 1))
$ws.Cells.Item(181, 1).Value = 179
$ws.Cells.Item(181, 2).Value = "Hellas Verona"
$ws.Cells.Item(181, 3).Value = "Udinese"
$ws.Cells.Item(181, 4).Value = 0
$ws.Cells.Item(181, 5).Value = 0
$ws.Cells.Item(181, 6).Value = 1.02
$ws.Cells.Item(181, 7).Value = 0.8100000000000001
$ws.Cells.Item(181, 8).Value = 1.2
$ws.Cells.Item(181, 9).Value = 0.93
$ws.Cells.Item(181, 10).Value = 0
$ws.Cells.Item(181, 11).Value = 0
$ws.Cells.Item(181, 12).Value = 0.18
$ws.Cells.Item(181, 13).Value = 0.12
$ws.Cells.Item(181, 14).Value = 0.3
$ws.Cells.Item(181, 15).Value = 0

# Row 182
$ws.Cells.Item(181, 1).Copy($ws.Cells.Item(182, 1))
$ws.Cells.Item(182, 1).Value = 180
$ws.Cells.Item(182, 2).Value = "Venezia"
$ws.Cells.Item(182, 3).Value = "Empoli"
$ws.Cells.Item(182, 4).Value = 1
$ws.Cells.Item(182, 5).Value = 1
$ws.Cells.Item(182, 6).Value = 0.83
$ws.Cells.Item(182, 7).Value = 1.78
$ws.Cells.Item(182, 8).Value = 0.63
$ws.Cells.Item(182, 9).Value = 1.51
$ws.Cells.Item(182, 10).Value = 0
$ws.Cells.Item(182, 11).Value = 0
$ws.Cells.Item(182, 12).Value = 0.2
$ws.Cells.Item(182, 13).Value = 0.27
$ws.Cells.Item(182, 14).Value = 0.47
$ws.Cells.Item(182, 15).Value = 2

# Row 183
$ws.Cells.Item(182, 1).Copy($ws.Cells.Item(183, 1))
$ws.Cells.Item(183, 1).Value = 181
$ws.Cells.Item(183, 2).Value = "Lecce"
$ws.Cells.Item(183, 3).Value = "Genoa"
$ws.Cells.Item(183, 4).Value = 0
$ws.Cells.Item(183, 5).Value = 0
$ws.Cells.Item(183, 6).Value = 0.5600000000000001
$ws.Cells.Item(183, 7).Value = 1.2
$ws.Cells.Item(183, 8).Value = 0.85
$ws.Cells.Item(183, 9).Value = 1
$ws.Cells.Item(183, 10).Value = 0
$ws.Cells.Item(183, 11).Value = 1
$ws.Cells.Item(183, 12).Value = 0.29
$ws.Cells.Item(183, 13).Value = 0.2
$ws.Cells.Item(183, 14).Value = 0.49
$ws.Cells.Item(183, 15).Value = 1

# Row 184
$ws.Cells.Item(183, 1).Copy($ws.Cells.Item(184, 1))
$ws.Cells.Item(184, 1).Value = 182
$ws.Cells.Item(184, 2).Value = "Monza"
$ws.Cells.Item(184, 3).Value = "Cagliari"
$ws.Cells.Item(184, 4).Value = 1
$ws.Cells.Item(184, 5).Value = 2
$ws.Cells.Item(184, 6).Value = 1.58
$ws.Cells.Item(184, 7).Value = 0.85
$ws.Cells.Item(184, 8).Value = 1.3
$ws.Cells.Item(184, 9).Value = 1.18
$ws.Cells.Item(184, 10).Value = 1
$ws.Cells.Item(184, 11).Value = 0
$ws.Cells.Item(184, 12).Value = 0.28
$ws.Cells.Item(184, 13).Value = 0.33
$ws.Cells.Item(184, 14).Value = 0.6
$ws.Cells.Item(184, 15).Value = 2

# Row 185
$ws.Cells.Item(184, 1).Copy($ws.Cells.Item(185, 1))
$ws.Cells.Item(185, 1).Value = 183
$ws.Cells.Item(185, 2).Value = "Roma"
$ws.Cells.Item(185, 3).Value = "Lazio"
$ws.Cells.Item(185, 4).Value = 2
$ws.Cells.Item(185, 5).Value = 0
$ws.Cells.Item(185, 6).Value = 1.18
$ws.Cells.Item(185, 7).Value = 0.87
$ws.Cells.Item(185, 8).Value = 1.05
$ws.Cells.Item(185, 9).Value = 1.16
$ws.Cells.Item(185, 10).Value = 0
$ws.Cells.Item(185, 11).Value = 0
$ws.Cells.Item(185, 12).Value = 0.13
$ws.Cells.Item(185, 13).Value = 0.29
$ws.Cells.Item(185, 14).Value = 0.42
$ws.Cells.Item(185, 15).Value = 2

# Row 186
$ws.Cells.Item(185, 1).Copy($ws.Cells.Item(186, 1))
$ws.Cells.Item(186, 1).Value = 184
$ws.Cells.Item(186, 2).Value = "Torino"
$ws.Cells.Item(186, 3).Value = "Parma"
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 0
$ws.Cells.Item(186, 6).Value = 1.08
$ws.Cells.Item(186, 7).Value = 0.38
$ws.Cells.Item(186, 8).Value = 1.47
$ws.Cells.Item(186, 9).Value = 0.38
$ws.Cells.Item(186, 10).Value = 0
$ws.Cells.Item(186, 11).Value = 0
$ws.Cells.Item(186, 12).Value = 0.39
$ws.Cells.Item(186, 13).Value = 0
$ws.Cells.Item(186, 14).Value = 0.39
$ws.Cells.Item(186, 15).Value = 0

# Row 187
$ws.Cells.Item(186, 1).Copy($ws.Cells.Item(187, 1))
$ws.Cells.Item(187, 1).Value = 185
$ws.Cells.Item(187, 2).Value = "Lazio"
$ws.Cells.Item(187, 3).Value = "Como"
$ws.Cells.Item(187, 4).Value = 1
$ws.Cells.Item(187, 5).Value = 1
$ws.Cells.Item(187, 6).Value = 0.53
$ws.Cells.Item(187, 7).Value = 1.98
$ws.Cells.Item(187, 8).Value = 0.77
$ws.Cells.Item(187, 9).Value = 1.62
$ws.Cells.Item(187, 10).Value = 0
$ws.Cells.Item(187, 11).Value = 0
$ws.Cells.Item(187, 12).Value = 0.24
$ws.Cells.Item(187, 13).Value = 0.36
$ws.Cells.Item(187, 14).Value = 0.6
$ws.Cells.Item(187, 15).Value = 2

# Row 188
$ws.Cells.Item(187, 1).Copy($ws.Cells.Item(188, 1))
$ws.Cells.Item(188, 1).Value = 186
$ws.Cells.Item(188, 2).Value = "Empoli"
$ws.Cells.Item(188, 3).Value = "Lecce"
$ws.Cells.Item(188, 4).Value = 1
$ws.Cells.Item(188, 5).Value = 3
$ws.Cells.Item(188, 6).Value = 1.43
$ws.Cells.Item(188, 7).Value = 1.46
$ws.Cells.Item(188, 8).Value = 1.61
$ws.Cells.Item(188, 9).Value = 1.23
$ws.Cells.Item(188, 10).Value = 0
$ws.Cells.Item(188, 11).Value = 0
$ws.Cells.Item(188, 12).Value = 0.18
$ws.Cells.Item(188, 13).Value = 0.23
$ws.Cells.Item(188, 14).Value = 0.41
$ws.Cells.Item(188, 15).Value = 4

# Row 189
$ws.Cells.Item(188, 1).Copy($ws.Cells.Item(189, 1))
$ws.Cells.Item(189, 1).Value = 187
$ws.Cells.Item(189, 2).Value = "Milan"
$ws.Cells.Item(189, 3).Value = "Cagliari"
$ws.Cells.Item(189, 4).Value = 1
$ws.Cells.Item(189, 5).Value = 1
$ws.Cells.Item(189, 6).Value = 2.85
$ws.Cells.Item(189, 7).Value = 0.43
$ws.Cells.Item(189, 8).Value = 2.58
$ws.Cells.Item(189, 9).Value = 0.4
$ws.Cells.Item(189, 10).Value = 1
$ws.Cells.Item(189, 11).Value = 0
$ws.Cells.Item(189, 12).Value = 0.27
$ws.Cells.Item(189, 13).Value = 0.03
$ws.Cells.Item(189, 14).Value = 0.3
$ws.Cells.Item(189, 15).Value = 1

# Row 190
$ws.Cells.Item(189, 1).Copy($ws.Cells.Item(190, 1))
$ws.Cells.Item(190, 1).Value = 188
$ws.Cells.Item(190, 2).Value = "Torino"
$ws.Cells.Item(190, 3).Value = "Juventus"
$ws.Cells.Item(190, 4).Value = 1
$ws.Cells.Item(190, 5).Value = 1
$ws.Cells.Item(190, 6).Value = 0.59
$ws.Cells.Item(190, 7).Value = 0.9
$ws.Cells.Item(190, 8).Value = 0.79
$ws.Cells.Item(190, 9).Value = 1.48
$ws.Cells.Item(190, 10).Value = 0
$ws.Cells.Item(190, 11).Value = 0
$ws.Cells.Item(190, 12).Value = 0.2
$ws.Cells.Item(190, 13).Value = 0.58
$ws.Cells.Item(190, 14).Value = 0.78
$ws.Cells.Item(190, 15).Value = 2

# Row 191
$ws.Cells.Item(190, 1).Copy($ws.Cells.Item(191, 1))
$ws.Cells.Item(191, 1).Value = 189
$ws.Cells.Item(191, 2).Value = "Udinese"
$ws.Cells.Item(191, 3).Value = "Atalanta"
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 0
$ws.Cells.Item(191, 6).Value = 1.83
$ws.Cells.Item(191, 7).Value = 0.23
$ws.Cells.Item(191, 8).Value = 1.75
$ws.Cells.Item(191, 9).Value = 0.29
$ws.Cells.Item(191, 10).Value = 1
$ws.Cells.Item(191, 11).Value = 0
$ws.Cells.Item(191, 12).Value = 0.08
$ws.Cells.Item(191, 13).Value = 0.06
$ws.Cells.Item(191, 14).Value = 0.14
$ws.Cells.Item(191, 15).Value = 1

# Row 192
$ws.Cells.Item(191, 1).Copy($ws.Cells.Item(192, 1))
$ws.Cells.Item(192, 1).Value = 190
$ws.Cells.Item(192, 2).Value = "Bologna"
$ws.Cells.Item(192, 3).Value = "Roma"
$ws.Cells.Item(192, 4).Value = 2
$ws.Cells.Item(192, 5).Value = 2
$ws.Cells.Item(192, 6).Value = 2.06
$ws.Cells.Item(192, 7).Value = 1.4
$ws.Cells.Item(192, 8).Value = 1.84
$ws.Cells.Item(192, 9).Value = 1.34
$ws.Cells.Item(192, 10).Value = 1
$ws.Cells.Item(192, 11).Value = 1
$ws.Cells.Item(192, 12).Value = 0.22
$ws.Cells.Item(192, 13).Value = 0.06
$ws.Cells.Item(192, 14).Value = 0.28
$ws.Cells.Item(192, 15).Value = 2

# Row 193
$ws.Cells.Item(192, 1).Copy($ws.Cells.Item(193, 1))
$ws.Cells.Item(193, 1).Value = 191
$ws.Cells.Item(193, 2).Value = "Genoa"
$ws.Cells.Item(193, 3).Value = "Parma"
$ws.Cells.Item(193, 4).Value = 1
$ws.Cells.Item(193, 5).Value = 0
$ws.Cells.Item(193, 6).Value = 1.19
$ws.Cells.Item(193, 7).Value = 0.46
$ws.Cells.Item(193, 8).Value = 1.06
$ws.Cells.Item(193, 9).Value = 0.62
$ws.Cells.Item(193, 10).Value = 0
$ws.Cells.Item(193, 11).Value = 0
$ws.Cells.Item(193, 12).Value = 0.13
$ws.Cells.Item(193, 13).Value = 0.16
$ws.Cells.Item(193, 14).Value = 0.29
$ws.Cells.Item(193, 15).Value = 1

# Row 194
$ws.Cells.Item(193, 1).Copy($ws.Cells.Item(194, 1))
$ws.Cells.Item(194, 1).Value = 192
$ws.Cells.Item(194, 2).Value = "Napoli"
$ws.Cells.Item(194, 3).Value = "Hellas Verona"
$ws.Cells.Item(194, 4).Value = 2
$ws.Cells.Item(194, 5).Value = 0
$ws.Cells.Item(194, 6).Value = 1.18
$ws.Cells.Item(194, 7).Value = 0.27
$ws.Cells.Item(194, 8).Value = 1.53
$ws.Cells.Item(194, 9).Value = 0.33
$ws.Cells.Item(194, 10).Value = 0
$ws.Cells.Item(194, 11).Value = 0
$ws.Cells.Item(194, 12).Value = 0.35
$ws.Cells.Item(194, 13).Value = 0.06
$ws.Cells.Item(194, 14).Value = 0.41
$ws.Cells.Item(194, 15).Value = 2

# Row 195
$ws.Cells.Item(194, 1).Copy($ws.Cells.Item(195, 1))
$ws.Cells.Item(195, 1).Value = 193
$ws.Cells.Item(195, 2).Value = "Venezia"
$ws.Cells.Item(195, 3).Value = "Inter"
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 5).Value = 1
$ws.Cells.Item(195, 6).Value = 0.41
$ws.Cells.Item(195, 7).Value = 1.93
$ws.Cells.Item(195, 8).Value = 0.59
$ws.Cells.Item(195, 9).Value = 1.58
$ws.Cells.Item(195, 10).Value = 0
$ws.Cells.Item(195, 11).Value = 0
$ws.Cells.Item(195, 12).Value = 0.18
$ws.Cells.Item(195, 13).Value = 0.35
$ws.Cells.Item(195, 14).Value = 0.53
$ws.Cells.Item(195, 15).Value = 1

Write-Output "Added rows 180-195 (matches up through round 20)"
